$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "287.86"
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "-0.89%"
$ws.Cells.Item(2, 7).NumberFormat = "@"
$ws.Cells.Item(2, 7).Value = "11"

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "30.96"
$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = "0.70%"
$ws.Cells.Item(3, 7).NumberFormat = "@"
$ws.Cells.Item(3, 7).Value = "11"

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "4.917"
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = "-0.71%"
$ws.Cells.Item(4, 7).NumberFormat = "@"
$ws.Cells.Item(4, 7).Value = "11"

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "0.07317"
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = "1.41%"
$ws.Cells.Item(5, 7).NumberFormat = "@"
$ws.Cells.Item(5, 7).Value = "11"

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "2.330"
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = "25.16%"
$ws.Cells.Item(6, 7).NumberFormat = "@"
$ws.Cells.Item(6, 7).Value = "11"

$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = "0.59%"
$ws.Cells.Item(7, 7).NumberFormat = "@"
$ws.Cells.Item(7, 7).Value = "11"

$ws.Cells.Item(8, 2).Value = "GateToken"
$ws.Cells.Item(8, 3).Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "3.722"
$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = "-1.29%"
$ws.Cells.Item(8, 7).NumberFormat = "@"
$ws.Cells.Item(8, 7).Value = "11"

$ws.Cells.Item(9, 2).Value = "MXToken"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.9033"
$ws.Cells.Item(9, 5).NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = "0.62%"
$ws.Cells.Item(9, 7).NumberFormat = "@"
$ws.Cells.Item(9, 7).Value = "11"

$ws.Cells.Item(10, 2).Value = "LiechtensteinCryptoassetsExchange"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.09142"
$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = "18.29%"
$ws.Cells.Item(10, 7).NumberFormat = "@"
$ws.Cells.Item(10, 7).Value = "11"

$ws.Cells.Item(11, 2).Value = "WazirX"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.1689"
$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = "1.56%"
$ws.Cells.Item(11, 7).NumberFormat = "@"
$ws.Cells.Item(11, 7).Value = "11"

$ws.Cells.Item(12, 2).Value = "MandalaExchangeToken"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.08183"
$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = "2.94%"
$ws.Cells.Item(12, 7).NumberFormat = "@"
$ws.Cells.Item(12, 7).Value = "11"

$ws.Cells.Item(13, 2).Value = "BitrueCoin"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.03124"
$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = "2.67%"
$ws.Cells.Item(13, 7).NumberFormat = "@"
$ws.Cells.Item(13, 7).Value = "11"

$ws.Cells.Item(14, 2).Value = "BitMartToken"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.09943"
$ws.Cells.Item(14, 5).NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = "-0.74%"
$ws.Cells.Item(14, 7).NumberFormat = "@"
$ws.Cells.Item(14, 7).Value = "11"

$ws.Cells.Item(15, 2).Value = "BitForexToken"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.001498"
$ws.Cells.Item(15, 5).NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = "0.06%"
$ws.Cells.Item(15, 7).NumberFormat = "@"
$ws.Cells.Item(15, 7).Value = "11"

$ws.Cells.Item(16, 2).Value = "TigerCash"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.005786"
$ws.Cells.Item(16, 5).NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = "0.54%"
$ws.Cells.Item(16, 7).NumberFormat = "@"
$ws.Cells.Item(16, 7).Value = "11"

$ws.Cells.Item(17, 2).Value = "LEO"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "3.496"
$ws.Cells.Item(17, 5).NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = "0.92%"
$ws.Cells.Item(17, 7).NumberFormat = "@"
$ws.Cells.Item(17, 7).Value = "11"

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "2.098"
$ws.Cells.Item(18, 5).NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = "0.86%"
$ws.Cells.Item(18, 7).NumberFormat = "@"
$ws.Cells.Item(18, 7).Value = "11"

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.3330"
$ws.Cells.Item(19, 5).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = "0.28%"
$ws.Cells.Item(19, 7).NumberFormat = "@"
$ws.Cells.Item(19, 7).Value = "11"

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.1294"
$ws.Cells.Item(20, 5).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = "1.16%"
$ws.Cells.Item(20, 7).NumberFormat = "@"
$ws.Cells.Item(20, 7).Value = "11"

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "4.228"
$ws.Cells.Item(21, 5).NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = "4.35%"
$ws.Cells.Item(21, 7).NumberFormat = "@"
$ws.Cells.Item(21, 7).Value = "11"

$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = "-12.14%"
$ws.Cells.Item(22, 7).NumberFormat = "@"
$ws.Cells.Item(22, 7).Value = "11"

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.04505"
$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = "0.22%"
$ws.Cells.Item(23, 7).NumberFormat = "@"
$ws.Cells.Item(23, 7).Value = "11"

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "0.001210"
$ws.Cells.Item(24, 5).NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = "-0.41%"
$ws.Cells.Item(24, 7).NumberFormat = "@"
$ws.Cells.Item(24, 7).Value = "11"

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "0.004163"
$ws.Cells.Item(25, 5).NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = "-10.16%"
$ws.Cells.Item(25, 7).NumberFormat = "@"
$ws.Cells.Item(25, 7).Value = "11"

$ws.Cells.Item(26, 5).NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = "3.85%"
$ws.Cells.Item(26, 7).NumberFormat = "@"
$ws.Cells.Item(26, 7).Value = "11"

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0.0003394"
$ws.Cells.Item(27, 7).NumberFormat = "@"
$ws.Cells.Item(27, 7).Value = "11"

$ws.Cells.Item(28, 7).NumberFormat = "@"
$ws.Cells.Item(28, 7).Value = "11"

$ws.Cells.Item(29, 7).NumberFormat = "@"
$ws.Cells.Item(29, 7).Value = "11"

$ws.Cells.Item(30, 7).NumberFormat = "@"
$ws.Cells.Item(30, 7).Value = "11"

$ws.Cells.Item(31, 7).NumberFormat = "@"
$ws.Cells.Item(31, 7).Value = "11"

$ws.Cells.Item(32, 7).NumberFormat = "@"
$ws.Cells.Item(32, 7).Value = "11"

$ws.Cells.Item(33, 7).NumberFormat = "@"
$ws.Cells.Item(33, 7).Value = "11"

$ws.Cells.Item(34, 7).NumberFormat = "@"
$ws.Cells.Item(34, 7).Value = "11"

$ws.Cells.Item(35, 7).NumberFormat = "@"
$ws.Cells.Item(35, 7).Value = "11"

$ws.Cells.Item(36, 7).NumberFormat = "@"
$ws.Cells.Item(36, 7).Value = "11"

$ws.Cells.Item(37, 7).NumberFormat = "@"
$ws.Cells.Item(37, 7).Value = "11"

$ws.Cells.Item(38, 7).NumberFormat = "@"
$ws.Cells.Item(38, 7).Value = "11"

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.01574"
$ws.Cells.Item(39, 5).NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = "0.09%"
$ws.Cells.Item(39, 7).NumberFormat = "@"
$ws.Cells.Item(39, 7).Value = "11"

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.04440"
$ws.Cells.Item(40, 5).NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = "1.09%"
$ws.Cells.Item(40, 7).NumberFormat = "@"
$ws.Cells.Item(40, 7).Value = "11"

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.007359"
$ws.Cells.Item(41, 5).NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = "0.95%"
$ws.Cells.Item(41, 7).NumberFormat = "@"
$ws.Cells.Item(41, 7).Value = "11"

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.009478"
$ws.Cells.Item(42, 5).NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = "-4.57%"
$ws.Cells.Item(42, 7).NumberFormat = "@"
$ws.Cells.Item(42, 7).Value = "11"

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.1328"
$ws.Cells.Item(43, 5).NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = "1.86%"
$ws.Cells.Item(43, 7).NumberFormat = "@"
$ws.Cells.Item(43, 7).Value = "11"

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.002230"
$ws.Cells.Item(44, 5).NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = "7.57%"
$ws.Cells.Item(44, 7).NumberFormat = "@"
$ws.Cells.Item(44, 7).Value = "11"

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.008985"
$ws.Cells.Item(45, 5).NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = "-4.50%"
$ws.Cells.Item(45, 7).NumberFormat = "@"
$ws.Cells.Item(45, 7).Value = "11"

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.00006115"
$ws.Cells.Item(46, 7).NumberFormat = "@"
$ws.Cells.Item(46, 7).Value = "11"

$ws.Cells.Item(47, 5).NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = "-0.14%"
$ws.Cells.Item(47, 7).NumberFormat = "@"
$ws.Cells.Item(47, 7).Value = "11"

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "2.458"
$ws.Cells.Item(48, 5).NumberFormat = "@"
$ws.Cells.Item(48, 5).Value = "6.44%"
$ws.Cells.Item(48, 7).NumberFormat = "@"
$ws.Cells.Item(48, 7).Value = "11"

$ws.Cells.Item(49, 7).NumberFormat = "@"
$ws.Cells.Item(49, 7).Value = "11"

$ws.Cells.Item(50, 5).NumberFormat = "@"
$ws.Cells.Item(50, 5).Value = "-0.14%"
$ws.Cells.Item(50, 7).NumberFormat = "@"
$ws.Cells.Item(50, 7).Value = "11"

$ws.Cells.Item(51, 5).NumberFormat = "@"
$ws.Cells.Item(51, 5).Value = "-0.14%"
$ws.Cells.Item(51, 7).NumberFormat = "@"
$ws.Cells.Item(51, 7).Value = "11"

Write-Host "Updated symbol list"